$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Achat a faire" : refresh the shopping list (vendors / refs / prices)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Achat a faire")

# Row 3 : Alimentation solaire / panneau solaire now bought from Mouser
$ws1.Range("E3").Value = "Mouser"

# Row 4 : convertisseur now bought from Mouser (was semageek)
$ws1.Range("E4").Value = "Mouser"

# Row 5 : switched to a 5Ah battery bought on ebay, with a real price
$ws1.Range("D5").Value = "batterie 5Ah"
$ws1.Range("E5").Value = "ebay"
$ws1.Range("F5").Value = 8
$ws1.Range("G5").Value = 1
$ws1.Range("H5").Formula = "=G5*F5"

# Row 8 : BME280 sensor now bought from Mouser (was semageek)
$ws1.Range("E8").Value = "Mouser"

# Row 9 : SI1145 sensor now bought from Mouser, real unit price (10 instead of 5)
$ws1.Range("E9").Value = "Mouser"
$ws1.Range("F9").Value = 10
$ws1.Range("H9").Formula = "=G9*F9"

# Row 10 : the sparkfun hyperlink moves down to row 25, clear it here
$ws1.Range("J10").ClearContents()

# Row 13 (sunfounder / shield ethernet) moves down to row 14, so clear row 13
$ws1.Range("B13").ClearContents()
$ws1.Range("D13").ClearContents()
$ws1.Range("E13").ClearContents()
$ws1.Range("F13").ClearContents()
$ws1.Range("G13").ClearContents()
$ws1.Range("H13").ClearContents()

# Row 14 : now holds the sunfounder / shield ethernet purchase (bought on amazon)
$ws1.Range("B14").Value = "sunfounder"
$ws1.Range("D14").Value = "shield ethernet"
$ws1.Range("E14").Value = "amazon"
$ws1.Range("F14").Value = 13
$ws1.Range("G14").Value = 1
$ws1.Range("H14").Formula = "=G14*F14"

# Row 15 : "fer a souder" (soldering iron), bought on amazon, with a real price
$ws1.Range("D15").Value = "fer à souder"
$ws1.Range("E15").Value = "amazon"
$ws1.Range("F15").Value = 20
$ws1.Range("G15").Value = 1
$ws1.Range("H15").Formula = "=G15*F15"

# Row 16 : "fer a souder" label moved up to row 15, clear leftover label
$ws1.Range("D16").ClearContents()

# Row 17 : new antenna bought from a chinese Aliexpress seller (qty 2)
$ws1.Range("B17").Value = "chinois"
$ws1.Range("D17").Value = "antenne"
$ws1.Range("E17").Value = "Aliexpress"
$ws1.Range("F17").Value = 4.01
$ws1.Range("G17").Value = 2
$ws1.Range("H17").Formula = "=G17*F17"

# Row 18 : matching SMA pcb connector, also from Aliexpress
$ws1.Range("D18").Value = "connecteur SMA pcb"
$ws1.Range("E18").Value = "Aliexpress"
$ws1.Range("F18").Value = 1.49
$ws1.Range("G18").Value = 1
$ws1.Range("H18").Formula = "=G18*F18"

# Row 25 : the sparkfun hyperlink text now lives next to this row instead of J10
$ws1.Range("J10").Copy($ws1.Range("J25"))
$ws1.Range("J25").Value = "https://www.sparkfun.com/products/15441"

$ws1.Range("F10").Select()

# ---------------------------------------------------------------------------
# Sheet "Pin des Arduino" : annotate RTC / LCD rows with the sensor reference
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Pin des Arduino")
$ws4.Range("C26").Value = "BME280"
$ws4.Range("C27").Value = "SI1145"

$ws4.Activate()
$ws4.Range("C26").Select()
